$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the old "U.S. CENSUS BUREAU, 2012" source line (old row 55),
# push everything below down by one.
$ws.Rows(55).Insert()

# The hyperlinked "http://www.census.gov/econ/islandareas/" line (now row 57) is removed...
$ws.Rows(57).Delete()

# ...and re-inserted as a new, plain (non-hyperlinked) row further down, right after the
# now-blank gap row (A57).
$ws.Rows(58).Insert()

# New blank source line.
$ws.Range("A55").Value = ""

# Re-added source URL text, no longer a hyperlink.
$ws.Range("A58").Value = "http://www.census.gov/econ/islandareas/"

# The old long citation text is replaced with a short "USCB" label.
$ws.Range("A62").Value = "USCB"

# Drop the now-stale hyperlink definition (it used to point at A56).
$ws.Hyperlinks.Delete()
